$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, shifting existing rows 129..159 down to 130..160
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new weekly data point
$ws.Range("A129").Value = 9
$ws.Range("B129").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C129").Value = "Metropolitana"
$ws.Range("D129").Value = 44511
$ws.Range("E129").Value = 13
$ws.Range("F129").Value = 100112030
$ws.Range("G129").Value = "Poroto granado"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Primera"
$ws.Range("J129").Value = 25
$ws.Range("K129").Value = 35000
$ws.Range("L129").Value = 38000
$ws.Range("M129").Value = 36440
$ws.Range("N129").Value = "$/malla 25 kilos"
$ws.Range("O129").Value = "Perú"
$ws.Range("P129").Value = 1458
$ws.Range("Q129").Value = 25
$ws.Range("R129").Value = "Hortaliza"
